# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps produced by a fresh report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 261afc06 row (row 3)
$wsOverview.Range("G3").Value = "2016-10-20 00:10:50"

# zh-cn sheet, 261afc06 row (row 3)
$wsZhCn.Range("H3").Value = "2016-10-20 00:10:38"   # Correspond Handoff Datetime
$wsZhCn.Range("K3").Value = "2016-10-20 00:11:27"   # Correspond Handback DateTime

# de-de sheet, 261afc06 row (row 3)
$wsDeDe.Range("H3").Value = "2016-10-20 00:10:50"   # Correspond Handoff Datetime
$wsDeDe.Range("K3").Value = "2016-10-20 00:11:47"   # Correspond Handback DateTime
